$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "LoginTestData"

# Set cell values
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "test@fastcollab.com"
$ws.Range("B2").Value = "Myadmin@123"

# Apply header style
$ws.Range("A1:B1").Style = "Accent5"

# Column widths (closest achievable values given engine's 1/6-character quantization)
$ws.Columns("A").ColumnWidth = 18.666666666666668
$ws.Columns("B").ColumnWidth = 20.333333333333336

# Selection
[void]$ws.Range("B6").Select()
